$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 23812574
$ws.Range("I100").Value = 33335554
$ws.Range("J100").Value = 5125
$ws.Range("K100").Value = 33335554
$ws.Range("L100").Value = 5125
$ws.Range("M100").Value = -33335013
$ws.Range("N100").Value = -6207
$ws.Range("H112").Value = 3967.3096
$ws.Range("J112").Value = 4056.756
$ws.Range("L112").Value = 12170.268
$ws.Range("N112").Value = -14386.268
$ws.Range("H137").Value = 903.5925999999999
$ws.Range("I137").Value = 846
$ws.Range("J137").Value = 1068.1428
$ws.Range("K137").Value = 2538
$ws.Range("L137").Value = 3204.4284
$ws.Range("M137").Value = 12
$ws.Range("N137").Value = -8304.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1377.7222
$ws.Range("I61").Value = 1127.9166
$ws.Range("J61").Value = 1877.3334
$ws.Range("K61").Value = 1127.9166
$ws.Range("L61").Value = 1877.3334
$ws.Range("M61").Value = -915.9166
$ws.Range("N61").Value = -2301.3334
$ws.Range("H74").Value = 1052.8572
$ws.Range("I74").Value = 1060.0741
$ws.Range("J74").Value = 1028.5
$ws.Range("K74").Value = 1060.0741
$ws.Range("L74").Value = 1028.5
$ws.Range("M74").Value = -186.0741
$ws.Range("N74").Value = -2776.5
$ws.Range("H77").Value = 1052.8572
$ws.Range("I77").Value = 1060.0741
$ws.Range("J77").Value = 1028.5
$ws.Range("K77").Value = 5300.3705
$ws.Range("L77").Value = 5142.5
$ws.Range("M77").Value = -932.3705
$ws.Range("N77").Value = -13878.5
$ws.Range("H132").Value = 24417016
$ws.Range("I132").Value = 50002064
$ws.Range("J132").Value = 50300.855
$ws.Range("K132").Value = 150006192
$ws.Range("L132").Value = 150902.565
$ws.Range("M132").Value = -150003662
$ws.Range("N132").Value = -155962.565
$ws.Range("H136").Value = 1377.7222
$ws.Range("I136").Value = 1127.9166
$ws.Range("J136").Value = 1877.3334
$ws.Range("K136").Value = 3383.7498
$ws.Range("L136").Value = 5632.0002
$ws.Range("M136").Value = -833.7498000000001
$ws.Range("N136").Value = -10732.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5129.3657
$ws.Range("I134").Value = 1789.3448
$ws.Range("J134").Value = 13201.083
$ws.Range("K134").Value = 5368.0344
$ws.Range("L134").Value = 39603.249
$ws.Range("M134").Value = -2833.0344
$ws.Range("N134").Value = -44673.249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5002931.5
$ws.Range("I31").Value = 6758215.5
$ws.Range("J31").Value = 7123.6924
$ws.Range("K31").Value = 6758215.5
$ws.Range("L31").Value = 7123.6924
$ws.Range("M31").Value = -6757920.5
$ws.Range("N31").Value = -7713.6924
$ws.Range("H34").Value = 5002931.5
$ws.Range("I34").Value = 6758215.5
$ws.Range("J34").Value = 7123.6924
$ws.Range("K34").Value = 6758215.5
$ws.Range("L34").Value = 7123.6924
$ws.Range("M34").Value = -6758013.5
$ws.Range("N34").Value = -7527.6924
$ws.Range("H58").Value = 1010.25
$ws.Range("I58").Value = 934.93335
$ws.Range("K58").Value = 934.93335
$ws.Range("M58").Value = -731.93335
$ws.Range("H132").Value = 31852.117
$ws.Range("I132").Value = 1485.8334
$ws.Range("J132").Value = 104731.2
$ws.Range("K132").Value = 4457.5002
$ws.Range("L132").Value = 314193.6
$ws.Range("M132").Value = -1927.5002
$ws.Range("N132").Value = -319253.6
$ws.Range("H134").Value = 1819.0358
$ws.Range("I134").Value = 1145.6666
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 3436.9998
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -901.9998000000001
$ws.Range("N134").Value = -65070
$ws.Range("H136").Value = 1010.25
$ws.Range("I136").Value = 934.93335
$ws.Range("K136").Value = 2804.80005
$ws.Range("M136").Value = -254.8000499999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 251.68889
$ws.Range("J107").Value = 254.11111
$ws.Range("L107").Value = 762.3333299999999
$ws.Range("N107").Value = -4602.333329999999
$ws.Range("H137").Value = 3207.6
$ws.Range("I137").Value = 1396
$ws.Range("J137").Value = 4113.4
$ws.Range("K137").Value = 4188
$ws.Range("L137").Value = 12340.2
$ws.Range("M137").Value = 912
$ws.Range("N137").Value = -22540.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1505.7858
$ws.Range("I113").Value = 1268.1
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1268.1
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 901.9000000000001
$ws.Range("N113").Value = -6440
$ws.Range("H132").Value = 251325.25
$ws.Range("I132").Value = 39413.73
$ws.Range("J132").Value = 557419.7
$ws.Range("K132").Value = 118241.19
$ws.Range("L132").Value = 1672259.1
$ws.Range("M132").Value = -115711.19
$ws.Range("N132").Value = -1677319.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 6431.4375
$ws.Range("I131").Value = 2000
$ws.Range("J131").Value = 6726.8667
$ws.Range("K131").Value = 2000
$ws.Range("L131").Value = 6726.8667
$ws.Range("N131").Value = -16806.8667
$ws.Range("M131").Value = 3040
$ws.Range("H132").Value = 29581.611
$ws.Range("I132").Value = 51806.1
$ws.Range("J132").Value = 1801
$ws.Range("K132").Value = 155418.3
$ws.Range("L132").Value = 5403
$ws.Range("M132").Value = -152888.3
$ws.Range("N132").Value = -10463
$ws.Range("H136").Value = 8213.440000000001
$ws.Range("I136").Value = 8018.316
$ws.Range("J136").Value = 8831.333000000001
$ws.Range("K136").Value = 24054.948
$ws.Range("L136").Value = 26493.999
$ws.Range("M136").Value = -21504.948
$ws.Range("N136").Value = -31593.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1542.52
$ws.Range("I122").Value = 1498.65
$ws.Range("J122").Value = 1718
$ws.Range("K122").Value = 4495.950000000001
$ws.Range("L122").Value = 5154
$ws.Range("M122").Value = -2045.950000000001
$ws.Range("N122").Value = -10054
$ws.Range("H132").Value = 63678636
$ws.Range("I132").Value = 98262344
$ws.Range("J132").Value = 2492076.2
$ws.Range("K132").Value = 294787032
$ws.Range("L132").Value = 7476228.600000001
$ws.Range("M132").Value = -294784502
$ws.Range("N132").Value = -7481288.600000001
$ws.Range("H136").Value = 29107.723
$ws.Range("I136").Value = 40772.72
$ws.Range("J136").Value = 2596.3635
$ws.Range("K136").Value = 122318.16
$ws.Range("L136").Value = 7789.0905
$ws.Range("M136").Value = -119768.16
$ws.Range("N136").Value = -12889.0905
